# ============================================================================
# Adds a "Player Info" sheet (before "ODI Batting"), adds an
# "ODI Batting Extra" sheet (after "ODI Bowling"), renames the
# MATCH_CARD_LINK columns to MATCH_CODE on both existing sheets and replaces
# the howstat URL values with the bare match-code numbers.
#
# NOTE: this COM host re-indexes sheet objects on structural changes, so any
# worksheet reference captured *before* a Worksheets.Add() call must be
# re-fetched (by name) afterwards rather than reused - otherwise it can end
# up silently pointing at the wrong sheet. To stay safe, every structural
# operation (Add) happens first, and each sheet reference used for cell
# edits below is (re-)fetched by name immediately before it is used.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Structural changes: insert the two new sheets in their final positions.
# ---------------------------------------------------------------------------
$odiBattingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odiBattingForInsert)
$playerInfo.Name = "Player Info"

$odiBowlingForInsert = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBowlingForInsert)
$extra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# 2. Populate "Player Info"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4419"
$playerInfo.Range("B2").Value = "James M Vince"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------------
# 3. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "3800"
    3  = "3911"
    4  = "3946"
    5  = "3948"
    6  = "3949"
    7  = "4175"
    8  = "4284"
    9  = "4287"
    10 = "4297"
    11 = "4300"
    12 = "4326"
    13 = "4331"
    14 = "4336"
    15 = "4426"
    16 = "4427"
    17 = "4428"
    18 = "4472"
    19 = "4473"
    20 = "4476"
    21 = "4660"
    22 = "4663"
    23 = "4666"
    24 = "4711"
    25 = "4713"
    26 = "4717"
}

foreach ($row in $battingCodes.Keys) {
    $cell = $odiBatting.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
}

# Rows where the player did not bat have always-empty INNING_NUMBER cells;
# re-assert that (assigning "" drops the cell entirely, matching the target).
$odiBatting.Range("B2").Value = ""
$odiBatting.Range("B9").Value = ""
$odiBatting.Range("B18").Value = ""

# ---------------------------------------------------------------------------
# 4. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4427"
    3 = "4428"
}

foreach ($row in $bowlingCodes.Keys) {
    $cell = $odiBowling.Range("B" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$row]
}

# ---------------------------------------------------------------------------
# 5. Populate "ODI Batting Extra"
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160
$extraHeader.Borders.LineStyle = 1

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH are text;
# only BATTING_POSITION is stored as a real number (and left General so it
# actually lands in the sheet as a numeric cell rather than a text one).
$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:C21").NumberFormat = "@"
$extra.Range("D2:D21").NumberFormat = "@"
$extra.Range("E2:E21").NumberFormat = "@"
$extra.Range("F2:F21").NumberFormat = "@"

# Rows with no recorded batting position still keep an (empty) cell in B -
# format just those as text first so the blank assignment below doesn't
# drop the cell entirely.
$blankPositionRows = @(2, 4, 7, 10, 16, 17, 18)
foreach ($row in $blankPositionRows) {
    $extra.Range("B" + $row).NumberFormat = "@"
}

$extraRows = @(
    @{ Row = 2;  Code = "4175"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 3;  Code = "4284"; Pos = 1;     N4 = "3";   N6 = "0";   Pct = "9.05%";  MoM = "NO"  }
    @{ Row = 4;  Code = "4287"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 5;  Code = "4297"; Pos = 2;     N4 = "6";   N6 = "0";   Pct = "12.61%"; MoM = "NO"  }
    @{ Row = 6;  Code = "4300"; Pos = 1;     N4 = "7";   N6 = "0";   Pct = "9.40%";  MoM = "NO"  }
    @{ Row = 7;  Code = "4326"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 8;  Code = "4331"; Pos = 1;     N4 = "2";   N6 = "0";   Pct = "6.60%";  MoM = "NO"  }
    @{ Row = 9;  Code = "4336"; Pos = 2;     N4 = "0";   N6 = "0";   Pct = $null;    MoM = "NO"  }
    @{ Row = 10; Code = "4426"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 11; Code = "4427"; Pos = 3;     N4 = "3";   N6 = "0";   Pct = "7.41%";  MoM = "NO"  }
    @{ Row = 12; Code = "4428"; Pos = 3;     N4 = "3";   N6 = "0";   Pct = "4.88%";  MoM = "NO"  }
    @{ Row = 13; Code = "4472"; Pos = 4;     N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 14; Code = "4473"; Pos = 4;     N4 = "8";   N6 = "0";   Pct = "22.67%"; MoM = "NO"  }
    @{ Row = 15; Code = "4476"; Pos = 4;     N4 = "11";  N6 = "0";   Pct = "30.72%"; MoM = "YES" }
    @{ Row = 16; Code = "4660"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 17; Code = "4663"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 18; Code = "4666"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    MoM = "NO"  }
    @{ Row = 19; Code = "4711"; Pos = 4;     N4 = "1";   N6 = "0";   Pct = "2.83%";  MoM = "NO"  }
    @{ Row = 20; Code = "4713"; Pos = 4;     N4 = "0";   N6 = "0";   Pct = "1.53%";  MoM = "NO"  }
    @{ Row = 21; Code = "4717"; Pos = 4;     N4 = "3";   N6 = "1";   Pct = "19.39%"; MoM = "NO"  }
)

foreach ($r in $extraRows) {
    $row = $r.Row
    $extra.Range("A" + $row).Value = $r.Code

    if ($null -ne $r.Pos) {
        $extra.Range("B" + $row).Value = $r.Pos
    } else {
        $extra.Range("B" + $row).Value = ""
    }

    $extra.Range("C" + $row).Value = $(if ($null -ne $r.N4) { $r.N4 } else { "" })
    $extra.Range("D" + $row).Value = $(if ($null -ne $r.N6) { $r.N6 } else { "" })
    $extra.Range("E" + $row).Value = $(if ($null -ne $r.Pct) { $r.Pct } else { "" })
    $extra.Range("F" + $row).Value = $r.MoM
}
